$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Size"
# ---------------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
$wsSize.Range("G2").Value = 1729.802117575254
$wsSize.Range("G3").Value = 6553.056633182304
$wsSize.Range("G4").Value = 98.11647660857044
$wsSize.Range("C5").Value = 297.8
$wsSize.Range("D5").Value = 1073.106
$wsSize.Range("E5").Value = 56.101
$wsSize.Range("F5").Value = 47.43703977891181
$wsSize.Range("G5").Value = 1474.444039778912

# ---------------------------------------------------------------------------
# Sheet "Cost"
# ---------------------------------------------------------------------------
$wsCost = $wb.Worksheets.Item("Cost")
$wsCost.Range("I2").Value = 5.675264476690022
$wsCost.Range("I3").Value = 0.3459604235150507
$wsCost.Range("I4").Value = 3.604181148250267
$wsCost.Range("I5").Value = 0.01962329532171408
$wsCost.Range("E6").Value = 0.02978
$wsCost.Range("F6").Value = 0.1073106
$wsCost.Range("G6").Value = 0.005610100000000001
$wsCost.Range("H6").Value = 0.004743703977891181
$wsCost.Range("I6").Value = 0.1474444039778912
$wsCost.Range("I7").Value = 0.09658249853902202
$wsCost.Range("I8").Value = 1.006186826078123
$wsCost.Range("I9").Value = 0.02739138298658309
$wsCost.Range("E10").Value = 0.006235308891553636
$wsCost.Range("F10").Value = 0.02246859430281919
$wsCost.Range("G10").Value = 0.001174637555826227
$wsCost.Range("H10").Value = 0.0009932323570263006
$wsCost.Range("I10").Value = 0.03087177310722536
$wsCost.Range("I11").Value = 0.3970227249141455

# ---------------------------------------------------------------------------
# Sheet "Indicators"
# ---------------------------------------------------------------------------
$wsInd = $wb.Worksheets.Item("Indicators")
$wsInd.Range("C2").Value = 296111.1311914226
$wsInd.Range("D2").Value = 145482.7919381901
$wsInd.Range("E2").Value = 441593.9231296126
$wsInd.Range("E3").Value = 0.3674766593047327
$wsInd.Range("E4").Value = 0.6325233406952673
$wsInd.Range("C5").Value = 0.1613387984024015
$wsInd.Range("D5").Value = 1.010097467358436
$wsInd.Range("E5").Value = 0.4409617152359359
$wsInd.Range("E6").Value = 0.0001252770048386342
